$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Item.getname unable to invoke"
$ws.Range("C7").Value = "ricky"
$ws.Range("D7").Value = "return value of location item getter was null"
$ws.Range("E7").Value = "added item to location"
$ws.Range("F7").Value = "fixed"

$ws.Range("F7").Select()
